$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -5
$ws.Range("F4").Value = 1
$ws.Range("F10").Value = -4
$ws.Range("F16").Value = -10
$ws.Range("F18").Value = 0
